$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.annual-entitlement"
$ws.Range("B3").Value2 = "Éves jogosultság"
$ws.Range("A4").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.annual-entitlement-changes"
$ws.Range("B4").Value2 = "Az alkalmazottak több szabadságot kérhetnek, mint amennyire jogosultak, és ha megadják őket, akkor negatív egyenlegük lesz. Ezt mindig figyelembe veszik a következő ciklus egyenlegének kiszámításakor."
$ws.Range("A5").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.annual-entitlement-value"
$ws.Range("B5").Value2 = ":óraszám"
$ws.Range("A6").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.carryover-label"
$ws.Range("B6").Value2 = "Átvitel"
$ws.Range("A7").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.carryover.expiry"
$ws.Range("B7").Value2 = "lejár :hónapokkal az időszak vége után"
$ws.Range("A8").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.carryover.expiry-none"
$ws.Range("B8").Value2 = "nincs lejárat"
$ws.Range("A9").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.carryover.limit"
$ws.Range("B9").Value2 = "Max. {{hours}}h"
$ws.Range("A10").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.carryover.not-allowed"
$ws.Range("B10").Value2 = "Nem alllowed"
$ws.Range("A11").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.granting-cycle-label"
$ws.Range("B11").Value2 = "Engedélyezési ciklus és pillanat"
$ws.Range("A12").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.granting-cycle.grant-at-end-of-the-month"
$ws.Range("B12").Value2 = "a következő hónap elején"
$ws.Range("A13").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.granting-cycle.grant-at-start-of-the-month"
$ws.Range("B13").Value2 = "a hónap elején"
$ws.Range("A14").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.granting-cycle.monthly"
$ws.Range("B14").Value2 = "Havi"
$ws.Range("A15").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.granting-cycle.yearly"
$ws.Range("B15").Value2 = "Éves"
$ws.Range("A16").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.policy-name"
$ws.Range("B16").Value2 = "Házirend neve"
$ws.Range("A17").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.proration-end-of-employment"
$ws.Range("B17").Value2 = "Proráció a foglalkoztatás végén"
$ws.Range("A18").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.proration-start-of-employment"
$ws.Range("B18").Value2 = "A foglalkoztatás megkezdésekor"
$ws.Range("A19").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.proration.daily"
$ws.Range("B19").Value2 = "Napi"
$ws.Range("A20").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.proration.monthly"
$ws.Range("B20").Value2 = "Havi"
$ws.Range("A21").Value2 = "absence-settings-ui.create-absence-policy-dialog.summary.proration.none"
$ws.Range("B21").Value2 = "Nincs"
$ws.Range("A22").Value2 = "absence-settings-ui.create-absence-policy-dialog.test.changes"
$ws.Range("B22").Value2 = "Korlátlan átvitel"
$ws.Range("A23").Value2 = "absence-settings-ui.create-absence-policy-test-changes"
$ws.Range("B23").Value2 = "Összefoglaló"
